# Update the "Time:" timestamp embedded in the statsmodels OLS summary text
# that is stored in cell B2 of each worksheet (sheets "46" down to "18").
# Sheet 1 ("46") gets the new time 20:59:35, all remaining sheets (2-29)
# get the new time 20:59:36 - matching the source diff.

$wb = $excel.ActiveWorkbook
$sheetCount = $wb.Worksheets.Count

for ($i = 1; $i -le $sheetCount; $i++) {
    $ws = $wb.Worksheets.Item($i)
    $cell = $ws.Range("B2")
    $text = $cell.Value()

    if ($null -eq $text) { continue }

    if ($i -eq 1) {
        $newText = $text.Replace("Time:                        20:51:34", "Time:                        20:59:35")
    } else {
        $newText = $text.Replace("Time:                        20:51:34", "Time:                        20:59:36")
        $newText = $newText.Replace("Time:                        20:51:35", "Time:                        20:59:36")
    }

    $cell.Value = $newText
}
